# The workbook's single data table (rows 2-33, columns A-T) gets its rows
# reshuffled: every data row's full contents moves to a (generally) different
# row position, while the header row (row 1) and the set of values used stay
# identical. This mapping was derived by diffing the original workbook
# against the target OOXML: $rowMap[$i] gives the ORIGINAL row number whose
# full contents now belongs at NEW row ($i + 2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 33
$firstCol = 1
$lastCol = 20

$rowMap = @(25, 5, 16, 8, 30, 23, 15, 28, 32, 11, 7, 22, 21, 27, 33, 18, 6, 3, 4, 19, 26, 20, 9, 17, 31, 13, 14, 24, 12, 2, 10, 29)

# Pass 1: snapshot every existing data row (by original row number) before
# any writes happen, since several destinations overlap with sources.
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowVals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals += , $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Pass 2: write each new row from the snapshot of its mapped original row.
for ($i = 0; $i -lt $rowMap.Length; $i++) {
    $newRow = $firstDataRow + $i
    $srcRow = $rowMap[$i]
    $vals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($newRow, $c).Value = $vals[$c - 1]
    }
}
